$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.432.82"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.849.16"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.76"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6298"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07671"
$ws.Range("E8").Value = "  +1.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2926"
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.69"
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07738"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "1.878.17"
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.026"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001073"
$ws.Range("E14").Value = "  +2.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6790"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.63"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "2.137.66"
$ws.Range("E17").Value = "  +2.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.191"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").Value = "29.459.79"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.98"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.435"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.81"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1379"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.343"
$ws.Range("E29").Value = "  +5.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.465"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05663"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.121"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.030"
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.843"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7087"
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.587"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01790"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("D40").Value = "1.220.28"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.558"
$ws.Range("E41").Value = "  +5.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9066"
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.72"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.21"
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("E46").Value = "  +1.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.132"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4021"
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.004"
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.680"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1145"
$ws.Range("E51").Value = "  +2.26%  "
